# Insert a new "Match ID" column before column A, shifting all existing
# columns one place to the right (A->B, B->C, ... W->X), then populate the
# new column with the constant Match ID value (16) for every data row and
# give it the bold "header" style used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right and leave column A blank.
$ws.Columns("A").Insert()

# Give the new column's header/data area the same bold font used by the
# other "Player ID"-style header/number cells (this creates the new bold,
# borderless cellXf that is applied below).
$ws.Range("A2:A18").Font.Bold = $true

# Row 2 is the header row - label the new column.
$ws.Range("A2").Value = "Match ID"

# Rows 4-18 hold the per-player (and totals) data; every one of them
# belongs to match 16.
for ($r = 4; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = 16
}

# Row 18 is the hidden "totals" row and keeps the default (non-bold) style
# for column A, matching the rest of that row.
$ws.Range("A18").Font.Bold = $false

# Touching cells on the hidden rows (1/3/18) can make Excel stamp an
# explicit custom row height; auto-fit them back so the rows keep their
# original (implicit) height.
$ws.Rows("3").AutoFit()
$ws.Rows("18").AutoFit()

# Match the saved selection state recorded in the target workbook.
[void]$ws.Range("A2:A17").Select()
